$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for Orégano was recorded. It belongs right
# after the existing row 189, so insert a blank row at 190 (this pushes the
# old rows 190-228 down to 191-229, carrying their formatting/styles along).
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with the new observation.
$ws.Cells.Item(190, 1).Value2 = 6
$ws.Cells.Item(190, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(190, 3).Value2 = "Metropolitana"
$ws.Cells.Item(190, 4).Value2 = 44798
$ws.Cells.Item(190, 5).Value2 = 13
$ws.Cells.Item(190, 6).Value2 = 100112029
$ws.Cells.Item(190, 7).Value2 = "Orégano"
$ws.Cells.Item(190, 8).Value2 = "Sin especificar"
$ws.Cells.Item(190, 9).Value2 = "Primera"
$ws.Cells.Item(190, 10).Value2 = 47
$ws.Cells.Item(190, 11).Value2 = 15000
$ws.Cells.Item(190, 12).Value2 = 16000
$ws.Cells.Item(190, 13).Value2 = 15468
$ws.Cells.Item(190, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(190, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(190, 16).Value2 = 5156
$ws.Cells.Item(190, 17).Value2 = 3
$ws.Cells.Item(190, 18).Value2 = "Hortaliza"
